$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): "Délai fournisseur" / "Multiple de quantité" were
# moved from the tail columns (F,G) to the front (C,D), pushing the
# quantity columns two slots to the right.
$ws.Range("C2").Value = "Multiple de quantité"
$ws.Range("D2").Value = "Délai fournisseur"
$ws.Range("E2").Value = "Quantité actuelle"
$ws.Range("F2").Value = "Quantité minimale"
$ws.Range("G2").Value = "Quantité maximale"

# --- Product reference code correction
$ws.Range("A4").Value = "FURN_8855"

# --- Column layout: widths were re-balanced to fit the new header text,
# and the two grouped/outlined columns (now holding the shorter
# "Délai fournisseur" / "Multiple de quantité" text at F:G) are shown
# instead of hidden.
$ws.Columns.Item(5).ColumnWidth = 16.666666666666668
$ws.Columns.Item(6).ColumnWidth = 14.833333333333334
$ws.Columns.Item(7).ColumnWidth = 18.0
$ws.Columns.Item(8).ColumnWidth = 11.0

$ws.Columns.Item(6).Hidden = $false
$ws.Columns.Item(7).Hidden = $false

# Columns K:L stay hidden (unchanged); re-assert explicitly since the
# engine only round-trips the hidden flag for columns touched in this
# session.
$ws.Columns.Item(11).Hidden = $true
$ws.Columns.Item(12).Hidden = $true
